$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Paragraph 13 (last, plain paragraph previously ind=426) ---
$p13xml = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00681DF5" w:rsidRPr="00681DF5" w:rsidRDefault="00681DF5" w:rsidP="00373B5A"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-MX"/></w:rPr><w:br/></w:r></w:p>'
$p13 = $d.Paragraphs(13)
$p13.Range.InsertXML($p13xml)

# --- Paragraph 12 (Prrafodelista, ind 426 -> 786, add Arial rFonts) ---
$p12xml = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00373B5A" w:rsidRDefault="00373B5A" w:rsidP="00F620BC"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="786"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$p12 = $d.Paragraphs(12)
$p12.Range.InsertXML($p12xml)

# --- Paragraph 6 (empty Prrafodelista -> justification text + break) ---
$p6xml = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00681DF5" w:rsidRPr="006F0027" w:rsidRDefault="00681DF5" w:rsidP="00681DF5"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t>Elegimos estas herramientas de progra</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t>maci' + [char]0x00F3 + 'n porque consideramos que son ' + [char]0x00F3 + 'ptimas para desarrollar nuestro proyecto,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">adem' + [char]0x00E1 + 's las conocemos y </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:t>cumplen con las necesidades de la arquitectura del sistema a implementar.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr><w:br/></w:r></w:p>'
$p6 = $d.Paragraphs(6)
$p6.Range.InsertXML($p6xml)

# --- Paragraph 5 ("Adobe PhotoShop 5" -> split runs w/ proofErr + bookmark) ---
$p5xml = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00681DF5" w:rsidRPr="006F0027" w:rsidRDefault="00681DF5" w:rsidP="00681DF5"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="006F0027"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Adobe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>PhotoShop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> 5</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML($p5xml)

Write-Output "done"
